# Updated cryptos list values (price + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.565.58"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "'1.813.37"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'305.44"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'0.4640"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").Value = "'0.3574"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "'46.26"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.8983"
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").Value = "'0.07767"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'19.33"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'1.831.08"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "'5.238"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'6.303"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'87.50"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'0.000008525"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'26.609.08"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "'14.15"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'4.978"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").Value = "'151.99"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'17.85"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'1.994"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "'112.90"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "'4.796"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'0.08727"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "'3.119"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").Value = "'2.753"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("D34").Value = "'0.7298"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").Value = "'4.418"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'1.121"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").Value = "'1.074"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'0.01927"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "'2.916"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").Value = "'0.05089"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "'0.5042"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("D42").Value = "'6.814"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'0.1491"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'7.950"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'0.4671"
$ws.Range("D46").Value = "'1.007"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'9.994"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Value = "'98.13"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'0.06020"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").Value = "'63.62"
$ws.Range("E51").Value = "  +0.56%  "
